$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.793150398273091
$ws.Range("D2").Value = 4.283966637138957
$ws.Range("E2").Value = 16.51581206488088
$ws.Range("F2").Value = 22.12710917420715
$ws.Range("G2").Value = 3.611507248062485
$ws.Range("K2").Value = 12.1898689693959
$ws.Range("O2").Value = 19.69776709310091

$ws.Range("B3").Value = 7.721051978771539
$ws.Range("D3").Value = 4.244083202001988
$ws.Range("E3").Value = 15.57734658471344
$ws.Range("F3").Value = 22.13571753165265
$ws.Range("G3").Value = 3.613781892252304
$ws.Range("K3").Value = 11.60425677748446
$ws.Range("O3").Value = 19.77170363363452

$ws.Range("B4").Value = 7.678304469731427
$ws.Range("D4").Value = 4.219243919358759
$ws.Range("E4").Value = 14.9762730581507
$ws.Range("F4").Value = 22.14977369286306
$ws.Range("G4").Value = 3.61525029970695
$ws.Range("K4").Value = 11.22720370185752
$ws.Range("O4").Value = 19.82380725407645

$ws.Range("B5").Value = 7.661286325873776
$ws.Range("D5").Value = 4.209039272153547
$ws.Range("E5").Value = 14.72535078157003
$ws.Range("F5").Value = 22.15769814008227
$ws.Range("G5").Value = 3.615866797302321
$ws.Range("K5").Value = 11.06927179968183
$ws.Range("O5").Value = 19.84671552071849

$ws.Range("B6").Value = 7.658485293552698
$ws.Range("D6").Value = 4.2073399882284
$ws.Range("E6").Value = 14.6833330728562
$ws.Range("F6").Value = 22.15914634628169
$ws.Range("G6").Value = 3.61597026189577
$ws.Range("K6").Value = 11.04279267223976
$ws.Range("O6").Value = 19.85062030282783

$ws.Range("B7").Value = 7.678073305230744
$ws.Range("D7").Value = 4.219106622097225
$ws.Range("E7").Value = 14.97291285025315
$ws.Range("F7").Value = 22.14987168486793
$ws.Range("G7").Value = 3.615258540600349
$ws.Range("K7").Value = 11.22509093559242
$ws.Range("O7").Value = 19.82410943301706

$ws.Range("B8").Value = 7.767986932106728
$ws.Range("D8").Value = 4.270291040836449
$ws.Range("E8").Value = 16.19753306065114
$ws.Range("F8").Value = 22.12825220610344
$ws.Range("G8").Value = 3.612276688394008
$ws.Range("K8").Value = 11.99164284852786
$ws.Range("O8").Value = 19.72186218762237

$ws.Range("B9").Value = 7.955496472215369
$ws.Range("D9").Value = 4.367633095611584
$ws.Range("E9").Value = 18.4852547266117
$ws.Range("F9").Value = 22.1557523535738
$ws.Range("G9").Value = 3.606995826406704
$ws.Range("K9").Value = 13.35194348402415
$ws.Range("O9").Value = 19.5750285022731

$ws.Range("B10").Value = 8.09886336303733
$ws.Range("D10").Value = 4.43697832447477
$ws.Range("E10").Value = 20.13044485292603
$ws.Range("F10").Value = 22.21884081654166
$ws.Range("G10").Value = 3.603457292364574
$ws.Range("K10").Value = 14.2597995382742
$ws.Range("O10").Value = 19.50047257446854

$ws.Range("B11").Value = 8.165035429096218
$ws.Range("D11").Value = 4.467986840393414
$ws.Range("E11").Value = 20.83671822943359
$ws.Range("F11").Value = 22.25686475224322
$ws.Range("G11").Value = 3.601920767349204
$ws.Range("K11").Value = 14.65225949558596
$ws.Range("O11").Value = 19.47390949278351

$ws.Range("B12").Value = 8.19020773127626
$ws.Range("D12").Value = 4.479646262455695
$ws.Range("E12").Value = 21.09814401971148
$ws.Range("F12").Value = 22.27260148347101
$ws.Range("G12").Value = 3.601349381174448
$ws.Range("K12").Value = 14.79788004104889
$ws.Range("O12").Value = 19.46491692200775

$ws.Range("B13").Value = 8.184781700509596
$ws.Range("D13").Value = 4.477138969010663
$ws.Range("E13").Value = 21.04210870690075
$ws.Range("F13").Value = 22.26915286611137
$ws.Range("G13").Value = 3.601471975032119
$ws.Range("K13").Value = 14.76665187800998
$ws.Range("O13").Value = 19.46680607317651

$ws.Range("B14").Value = 8.167104207463264
$ws.Range("D14").Value = 4.468947759737177
$ws.Range("E14").Value = 20.85834643527608
$ws.Range("F14").Value = 22.25813263465861
$ws.Range("G14").Value = 3.60187354972077
$ws.Range("K14").Value = 14.66430003838627
$ws.Range("O14").Value = 19.47314824932983

$ws.Range("B15").Value = 8.156290452983896
$ws.Range("D15").Value = 4.46391945650498
$ws.Range("E15").Value = 20.74500315829879
$ws.Range("F15").Value = 22.2515565291936
$ws.Range("G15").Value = 3.602120886578237
$ws.Range("K15").Value = 14.60121533445229
$ws.Range("O15").Value = 19.47717213365593

$ws.Range("B16").Value = 8.094556107310462
$ws.Range("D16").Value = 4.43494055658422
$ws.Range("E16").Value = 20.0834433290267
$ws.Range("F16").Value = 22.21654329752543
$ws.Range("G16").Value = 3.603559175197773
$ws.Range("K16").Value = 14.23373434263626
$ws.Range("O16").Value = 19.5023572719758

$ws.Range("B17").Value = 8.056912154227996
$ws.Range("D17").Value = 4.417021368777747
$ws.Range("E17").Value = 19.66683101266519
$ws.Range("F17").Value = 22.19745070715419
$ws.Range("G17").Value = 3.604460217300669
$ws.Range("K17").Value = 14.00300492874529
$ws.Range("O17").Value = 19.51969733924661

$ws.Range("B18").Value = 8.03535148391078
$ws.Range("D18").Value = 4.40666451014402
$ws.Range("E18").Value = 19.42323896469269
$ws.Range("F18").Value = 22.18734716437935
$ws.Range("G18").Value = 3.604985363706231
$ws.Range("K18").Value = 13.86836522472677
$ws.Range("O18").Value = 19.53036246705981

$ws.Range("B19").Value = 8.028067758421994
$ws.Range("D19").Value = 4.403149404549535
$ws.Range("E19").Value = 19.34008039606956
$ws.Range("F19").Value = 22.18407713020196
$ws.Range("G19").Value = 3.605164354580164
$ws.Range("K19").Value = 13.822448490085
$ws.Range("O19").Value = 19.53409199232285

$ws.Range("B20").Value = 8.060910150174172
$ws.Range("D20").Value = 4.418934135147139
$ws.Range("E20").Value = 19.71159049698633
$ws.Range("F20").Value = 22.19939228933958
$ws.Range("G20").Value = 3.60436358708805
$ws.Range("K20").Value = 14.02776658630608
$ws.Range("O20").Value = 19.51777981655941

$ws.Range("B21").Value = 8.172293592677798
$ws.Range("D21").Value = 4.47135600845392
$ws.Range("E21").Value = 20.91248510599658
$ws.Range("F21").Value = 22.26133326914454
$ws.Range("G21").Value = 3.601755314024385
$ws.Range("K21").Value = 14.69444485575765
$ws.Range("O21").Value = 19.47125639119005

$ws.Range("B22").Value = 8.245745297028327
$ws.Range("D22").Value = 4.505131437593392
$ws.Range("E22").Value = 21.66223869682883
$ws.Range("F22").Value = 22.30961050689012
$ws.Range("G22").Value = 3.600111612850251
$ws.Range("K22").Value = 15.11268320390847
$ws.Range("O22").Value = 19.44706952315168

$ws.Range("B23").Value = 8.206490199151812
$ws.Range("D23").Value = 4.487151076472728
$ws.Range("E23").Value = 21.26528164505243
$ws.Range("F23").Value = 22.28313236243122
$ws.Range("G23").Value = 3.600983329208055
$ws.Range("K23").Value = 14.89107260654255
$ws.Range("O23").Value = 19.45940668062938

$ws.Range("B24").Value = 8.059102401115704
$ws.Range("D24").Value = 4.418069543944076
$ws.Range("E24").Value = 19.69136743354468
$ws.Range("F24").Value = 22.19851178024275
$ws.Range("G24").Value = 3.604407251406613
$ws.Range("K24").Value = 14.01657803696699
$ws.Range("O24").Value = 19.51864456083441

$ws.Range("B25").Value = 7.903696415035684
$ws.Range("D25").Value = 4.341657799802075
$ws.Range("E25").Value = 17.84181157643668
$ws.Range("F25").Value = 22.14079293537749
$ws.Range("G25").Value = 3.608364208569877
$ws.Range("K25").Value = 12.99975903782674
$ws.Range("O25").Value = 19.60894074031387
